$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 562, shifting existing rows 562-574 down to 564-576
$ws.Rows.Item(562).Insert()
$ws.Rows.Item(563).Insert()

# New row 562 data
$ws.Range("A562").Value = 8
$ws.Range("B562").Value = "Terminal La Palmera de La Serena"
$ws.Range("C562").Value = "Coquimbo"
$ws.Range("D562").Value = 44448
$ws.Range("E562").Value = 4
$ws.Range("F562").Value = 100112020
$ws.Range("G562").Value = "Tomate"
$ws.Range("H562").Value = "Larga vida"
$ws.Range("I562").Value = "Primera"
$ws.Range("J562").Value = 860
$ws.Range("K562").Value = 8800
$ws.Range("L562").Value = 9000
$ws.Range("M562").Value = 8900
$ws.Range("N562").Value = "$/caja 10 kilos"
$ws.Range("O562").Value = "Región de Arica y Parinacota"
$ws.Range("P562").Value = 890
$ws.Range("Q562").Value = 10
$ws.Range("R562").Value = "Hortaliza"

# New row 563 data
$ws.Range("A563").Value = 8
$ws.Range("B563").Value = "Terminal La Palmera de La Serena"
$ws.Range("C563").Value = "Coquimbo"
$ws.Range("D563").Value = 44448
$ws.Range("E563").Value = 4
$ws.Range("F563").Value = 100112020
$ws.Range("G563").Value = "Tomate"
$ws.Range("H563").Value = "Larga vida"
$ws.Range("I563").Value = "Segunda"
$ws.Range("J563").Value = 520
$ws.Range("K563").Value = 7800
$ws.Range("L563").Value = 8000
$ws.Range("M563").Value = 7900
$ws.Range("N563").Value = "$/caja 10 kilos"
$ws.Range("O563").Value = "Región de Arica y Parinacota"
$ws.Range("P563").Value = 790
$ws.Range("Q563").Value = 10
$ws.Range("R563").Value = "Hortaliza"

# Apply the date number format (same as other date cells in column D) to the two new date cells
$ws.Range("D562:D563").NumberFormat = $ws.Range("D564").NumberFormat
